$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 750
$ws.Range("J29").Value = 1000
$ws.Range("L29").Value = 3000
$ws.Range("N29").Value = -3562
$ws.Range("H31").Value = 1150.25
$ws.Range("I31").Value = 865.6667
$ws.Range("J31").Value = 2004
$ws.Range("K31").Value = 2597.0001
$ws.Range("L31").Value = 6012
$ws.Range("M31").Value = -2367.0001
$ws.Range("N31").Value = -6472
$ws.Range("H38").Value = 542
$ws.Range("I38").Value = 710.8889
$ws.Range("K38").Value = 2132.6667
$ws.Range("M38").Value = -1760.6667
$ws.Range("H121").Value = 1115.1538
$ws.Range("J121").Value = 1115.1538
$ws.Range("L121").Value = 3345.4614
$ws.Range("N121").Value = -6839.4614
$ws.Range("H137").Value = 64900.5
$ws.Range("I137").Value = 1580
$ws.Range("J137").Value = 79512.92
$ws.Range("K137").Value = 4740
$ws.Range("L137").Value = 238538.76
$ws.Range("M137").Value = -2190
$ws.Range("N137").Value = -243638.76
$ws.Range("H138").Value = 2662.0645
$ws.Range("J138").Value = 2492.5
$ws.Range("L138").Value = 7477.5
$ws.Range("N138").Value = -17757.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2698.0781
$ws.Range("I32").Value = 1903.5962
$ws.Range("K32").Value = 1903.5962
$ws.Range("M32").Value = -1616.5962
$ws.Range("H37").Value = 13000
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents()
$ws.Range("H45").Value = 3129.8
$ws.Range("I45").Value = 5000
$ws.Range("J45").Value = 2799.7646
$ws.Range("K45").Value = 5000
$ws.Range("L45").Value = 2799.7646
$ws.Range("M45").Value = -4623
$ws.Range("N45").Value = -3553.7646
$ws.Range("H58").Value = 20000
$ws.Range("I58").Value = 20000
$ws.Range("K58").Value = 20000
$ws.Range("M58").Value = -19570
$ws.Range("H74").Value = 2342.2222
$ws.Range("I74").Value = 820
$ws.Range("J74").Value = 3560
$ws.Range("K74").Value = 820
$ws.Range("L74").Value = 3560
$ws.Range("M74").Value = 54
$ws.Range("N74").Value = -5308
$ws.Range("H77").Value = 2342.2222
$ws.Range("I77").Value = 820
$ws.Range("J77").Value = 3560
$ws.Range("K77").Value = 4100
$ws.Range("L77").Value = 17800
$ws.Range("M77").Value = 268
$ws.Range("N77").Value = -26536
$ws.Range("H80").Value = 49000
$ws.Range("J80").Value = 49000
$ws.Range("L80").Value = 49000
$ws.Range("N80").Value = -50996
$ws.Range("H83").Value = 49000
$ws.Range("J83").Value = 49000
$ws.Range("L83").Value = 147000
$ws.Range("N83").Value = -156984

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 3000
$ws.Range("I35").Value = 3000
$ws.Range("K35").Value = 3000
$ws.Range("M35").Value = -2690
$ws.Range("H134").Value = 3248
$ws.Range("I134").Value = 2897
$ws.Range("K134").Value = 8691
$ws.Range("M134").Value = -6156

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1411.1515
$ws.Range("I31").Value = 877.7
$ws.Range("J31").Value = 2231.8462
$ws.Range("K31").Value = 877.7
$ws.Range("L31").Value = 2231.8462
$ws.Range("M31").Value = -582.7
$ws.Range("N31").Value = -2821.8462
$ws.Range("H34").Value = 1411.1515
$ws.Range("I34").Value = 877.7
$ws.Range("J34").Value = 2231.8462
$ws.Range("K34").Value = 877.7
$ws.Range("L34").Value = 2231.8462
$ws.Range("M34").Value = -675.7
$ws.Range("N34").Value = -2635.8462
$ws.Range("H107").Value = 550.875
$ws.Range("I107").Value = 468.63635
$ws.Range("J107").Value = 731.8
$ws.Range("K107").Value = 468.63635
$ws.Range("L107").Value = 731.8
$ws.Range("M107").Value = 1451.36365
$ws.Range("N107").Value = -4571.8
$ws.Range("H132").Value = 2940.2
$ws.Range("I132").Value = 2678.7144
$ws.Range("J132").Value = 3273
$ws.Range("K132").Value = 8036.1432
$ws.Range("L132").Value = 9819
$ws.Range("M132").Value = -5506.1432
$ws.Range("N132").Value = -14879

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 5000
$ws.Range("I55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("M55").ClearContents()
$ws.Range("H92").Value = 994.8570999999999
$ws.Range("J92").Value = 994.8570999999999
$ws.Range("L92").Value = 2984.5713
$ws.Range("N92").Value = -5480.5713
$ws.Range("H97").Value = 2068.3333
$ws.Range("I97").Value = 1002
$ws.Range("K97").Value = 3006
$ws.Range("M97").Value = -2510
$ws.Range("H137").Value = 5755.778
$ws.Range("I137").Value = 2310
$ws.Range("J137").Value = 6186.5
$ws.Range("K137").Value = 6930
$ws.Range("L137").Value = 18559.5
$ws.Range("M137").Value = -1830
$ws.Range("N137").Value = -28759.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2750
$ws.Range("I122").Value = 1000
$ws.Range("J122").Value = 4500
$ws.Range("K122").Value = 3000
$ws.Range("L122").Value = 13500
$ws.Range("M122").Value = -550
$ws.Range("N122").Value = -18400

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4492.5
$ws.Range("I7").Value = 4266
$ws.Range("K7").Value = 4266
$ws.Range("M7").Value = -4154
$ws.Range("H16").Value = 4208.316
$ws.Range("I16").Value = 4518.1177
$ws.Range("J16").Value = 1575
$ws.Range("K16").Value = 4518.1177
$ws.Range("L16").Value = 1575
$ws.Range("M16").Value = -4348.1177
$ws.Range("N16").Value = -1915
$ws.Range("H22").Value = 5996
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 5996
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 5996
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -6586
$ws.Range("H27").Value = 5996
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 5996
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 5996
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -6210
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()
$ws.Range("H82").Value = 3793.3333
$ws.Range("I82").Value = 390
$ws.Range("J82").Value = 5495
$ws.Range("K82").Value = 390
$ws.Range("L82").Value = 5495
$ws.Range("M82").Value = -29
$ws.Range("N82").Value = -6217
$ws.Range("H85").Value = 3793.3333
$ws.Range("I85").Value = 390
$ws.Range("J85").Value = 5495
$ws.Range("K85").Value = 390
$ws.Range("L85").Value = 5495
$ws.Range("M85").Value = 858
$ws.Range("N85").Value = -7991
$ws.Range("H93").Value = 699.5
$ws.Range("I93").Value = 699.5
$ws.Range("K93").Value = 699.5
$ws.Range("M93").Value = 548.5
$ws.Range("H126").Value = 4492.5
$ws.Range("I126").Value = 4266
$ws.Range("K126").Value = 12798
$ws.Range("M126").Value = -10328
$ws.Range("H130").Value = 30429
$ws.Range("J130").Value = 30429
$ws.Range("L130").Value = 30429
$ws.Range("N130").Value = -40469

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()
$ws.Range("H132").Value = 3659.2
$ws.Range("I132").Value = 3011.25
$ws.Range("J132").Value = 4399.7144
$ws.Range("K132").Value = 9033.75
$ws.Range("L132").Value = 13199.1432
$ws.Range("M132").Value = -6503.75
$ws.Range("N132").Value = -18259.1432
